# "update field beauty - MDLWL"
#
# Genetics.xlsx keeps a single data record on Sheet1 (A2:K2) whose "id"
# column (A) is a generated CA-xxxxxxxx token. This commit refreshes that
# id to a new token, leaving every other field/style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CA-C4ENCTM5"
